$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 1500
$ws.Range("I21").Value = 1500
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -1032
$ws.Range("H23").Value = 1500
$ws.Range("I23").Value = 1500
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1266
$ws.Range("H40").Value = 1739.1852
$ws.Range("I40").Value = 2029.9166
$ws.Range("J40").Value = 1506.6
$ws.Range("K40").Value = 2029.9166
$ws.Range("L40").Value = 1506.6
$ws.Range("M40").Value = -1854.9166
$ws.Range("N40").Value = -1856.6
$ws.Range("H43").Value = 233333780
$ws.Range("I43").Value = 333333730
$ws.Range("J43").Value = 83333840
$ws.Range("K43").Value = 333333730
$ws.Range("L43").Value = 83333840
$ws.Range("M43").Value = -333333661
$ws.Range("N43").Value = -83333978
$ws.Range("H64").Value = 2924.6
$ws.Range("I64").Value = 2880.75
$ws.Range("K64").Value = 2880.75
$ws.Range("M64").Value = -2632.75
$ws.Range("H67").Value = 2924.6
$ws.Range("I67").Value = 2880.75
$ws.Range("K67").Value = 2880.75
$ws.Range("M67").Value = -2022.75
$ws.Range("H112").Value = 1655.8334
$ws.Range("J112").Value = 1786.95
$ws.Range("L112").Value = 5360.85
$ws.Range("N112").Value = -7576.85
$ws.Range("H116").Value = 18000
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H137").Value = 23257346
$ws.Range("I137").Value = 1127.1389
$ws.Range("J137").Value = 142860750
$ws.Range("K137").Value = 3381.4167
$ws.Range("L137").Value = 428582250
$ws.Range("M137").Value = -831.4166999999998
$ws.Range("N137").Value = -428587350
$ws.Range("H138").Value = 3006.5566
$ws.Range("I138").Value = 2665.111
$ws.Range("J138").Value = 3138.257
$ws.Range("K138").Value = 7995.333
$ws.Range("L138").Value = 9414.771000000001
$ws.Range("M138").Value = -2855.333
$ws.Range("N138").Value = -19694.771
$ws.Range("H141").Value = 5090.6665
$ws.Range("I141").Value = 2598.7778
$ws.Range("J141").Value = 7582.5557
$ws.Range("K141").Value = 7796.3334
$ws.Range("L141").Value = 22747.6671
$ws.Range("M141").Value = -2616.3334
$ws.Range("N141").Value = -33107.6671
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33606.598
$ws.Range("I32").Value = 29971.152
$ws.Range("J32").Value = 41569.953
$ws.Range("K32").Value = 29971.152
$ws.Range("L32").Value = 41569.953
$ws.Range("M32").Value = -29684.152
$ws.Range("N32").Value = -42143.953
$ws.Range("H61").Value = 1874.3513
$ws.Range("I61").Value = 1496.5
$ws.Range("J61").Value = 2767.4546
$ws.Range("K61").Value = 1496.5
$ws.Range("L61").Value = 2767.4546
$ws.Range("M61").Value = -1284.5
$ws.Range("N61").Value = -3191.4546
$ws.Range("H136").Value = 1874.3513
$ws.Range("I136").Value = 1496.5
$ws.Range("J136").Value = 2767.4546
$ws.Range("K136").Value = 4489.5
$ws.Range("L136").Value = 8302.363799999999
$ws.Range("M136").Value = -1939.5
$ws.Range("N136").Value = -13402.3638
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2159.9546
$ws.Range("I86").Value = 2093.125
$ws.Range("K86").Value = 2093.125
$ws.Range("M86").Value = -970.125
$ws.Range("H89").Value = 2159.9546
$ws.Range("I89").Value = 2093.125
$ws.Range("K89").Value = 10465.625
$ws.Range("M89").Value = -4849.625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15877052
$ws.Range("I31").Value = 2281.3215
$ws.Range("J31").Value = 28576868
$ws.Range("K31").Value = 2281.3215
$ws.Range("L31").Value = 28576868
$ws.Range("M31").Value = -1986.3215
$ws.Range("N31").Value = -28577458
$ws.Range("H34").Value = 15877052
$ws.Range("I34").Value = 2281.3215
$ws.Range("J34").Value = 28576868
$ws.Range("K34").Value = 2281.3215
$ws.Range("L34").Value = 28576868
$ws.Range("M34").Value = -2079.3215
$ws.Range("N34").Value = -28577272
$ws.Range("H58").Value = 6411709
$ws.Range("I58").Value = 1446.9412
$ws.Range("K58").Value = 1446.9412
$ws.Range("M58").Value = -1243.9412
$ws.Range("H129").Value = 49749.5
$ws.Range("I129").Value = 49001
$ws.Range("J129").Value = 49999
$ws.Range("K129").Value = 49001
$ws.Range("L129").Value = 49999
$ws.Range("M129").Value = -44001
$ws.Range("N129").Value = -59999
$ws.Range("H136").Value = 6411709
$ws.Range("I136").Value = 1446.9412
$ws.Range("K136").Value = 4340.8236
$ws.Range("M136").Value = -1790.8236
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2530
$ws.Range("J55").Value = 2530
$ws.Range("L55").Value = 7590
$ws.Range("N55").Value = -7944
$ws.Range("H131").Value = 765.9299999999999
$ws.Range("I131").Value = 415
$ws.Range("J131").Value = 788.3298
$ws.Range("K131").Value = 1245
$ws.Range("L131").Value = 2364.9894
$ws.Range("M131").Value = 3795
$ws.Range("N131").Value = -12444.9894
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 20834494
$ws.Range("I113").Value = 50000804
$ws.Range("K113").Value = 50000804
$ws.Range("M113").Value = -49998634
$ws.Range("H129").Value = 41499.75
$ws.Range("J129").Value = 41499.75
$ws.Range("L129").Value = 41499.75
$ws.Range("N129").Value = -51499.75
$ws.Range("H132").Value = 4412.2856
$ws.Range("I132").Value = 4573.6763
$ws.Range("J132").Value = 3726.375
$ws.Range("K132").Value = 13721.0289
$ws.Range("L132").Value = 11179.125
$ws.Range("M132").Value = -11191.0289
$ws.Range("N132").Value = -16239.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 45456612
$ws.Range("I7").Value = 2091.875
$ws.Range("J7").Value = 166668670
$ws.Range("K7").Value = 2091.875
$ws.Range("L7").Value = 166668670
$ws.Range("M7").Value = -1979.875
$ws.Range("N7").Value = -166668894
$ws.Range("H46").Value = 1307.6428
$ws.Range("I46").Value = 766.6667
$ws.Range("J46").Value = 1372.56
$ws.Range("K46").Value = 766.6667
$ws.Range("L46").Value = 1372.56
$ws.Range("M46").Value = -578.6667
$ws.Range("N46").Value = -1748.56
$ws.Range("H82").Value = 2271.3157
$ws.Range("I82").Value = 2143.4666
$ws.Range("K82").Value = 2143.4666
$ws.Range("M82").Value = -1782.4666
$ws.Range("H85").Value = 2271.3157
$ws.Range("I85").Value = 2143.4666
$ws.Range("K85").Value = 2143.4666
$ws.Range("M85").Value = -895.4666000000002
$ws.Range("H126").Value = 45456612
$ws.Range("I126").Value = 2091.875
$ws.Range("J126").Value = 166668670
$ws.Range("K126").Value = 6275.625
$ws.Range("L126").Value = 500006010
$ws.Range("M126").Value = -3805.625
$ws.Range("N126").Value = -500010950
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 43486496
$ws.Range("I126").Value = 62509876
$ws.Range("J126").Value = 4480.5713
$ws.Range("K126").Value = 187529628
$ws.Range("L126").Value = 13441.7139
$ws.Range("M126").Value = -187527158
$ws.Range("N126").Value = -18381.7139

Write-Host "Applied all cell updates."
